$d = $word.ActiveDocument

# 1. Insert the new paragraph right after the "Nedan presenteras..." paragraph
#    (3rd paragraph of the document) with the relocated "Vi förväntar..." text.
$introPara = $d.Paragraphs(3).Range
$introPara.InsertParagraphAfter()
$newPara = $d.Paragraphs(4).Range
$newPara.Text = "Vi förväntar oss att ni återkommer med ett skriftligt svar på vårt klagomål och även beskriver vilka korrigerande åtgärder ni satt in för att rätta till identifierade brister i er efterlevnad av den svenska FSC standarden."

# 2. Remove the trailing two empty paragraphs and the now-duplicated
#    "Vi förväntar..." paragraph at the end of the document.
$lastCount = $d.Paragraphs.Count
$startPara = $d.Paragraphs($lastCount - 2)
$endPara = $d.Paragraphs($lastCount)
$trailRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
$trailRange.Delete()

# 3. Update the date shown in the "first page" header from 2023-11-13 to 2023-11-14.
$sec = $d.Sections(1)
$firstPageHeader = $sec.Headers(2)
$firstPageHeader.Range.Find.Execute("2023-11-13", $true, $false, $false, $false, $false, $true, 1, $false, "2023-11-14", 2)
